# "unify functions into single_run.m + finish roy's review"
#
# output_data/1_midi.xlsx is the per-run MIDI event table for run_num = 1
# (headers: run_num, block_num, time_stamp, note, is_on, ipi). Now that
# single_run.m computes time_stamp/note/ipi itself, the cached
# time_stamp (C), note (D) and ipi (F) columns written by the old,
# per-run export path are stale and are reset to 0 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 is the header row; data runs from row 2 to the last used row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

if ($lastRow -ge 2) {
    $ws.Range("C2:C$lastRow").Value = 0
    $ws.Range("D2:D$lastRow").Value = 0
    $ws.Range("F2:F$lastRow").Value = 0
}

Write-Output "Zeroed time_stamp/note/ipi for rows 2..$lastRow on $($ws.Name)"
